$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Periodo Mora" column (E) for rows 17-25 previously listed periods
# 2206..2302 in ascending order. The workbook's records were refreshed:
# previous account-statement periods were removed and new ones added,
# which results in the periods now being listed in descending order
# (most recent period first).
$ws.Range("E17").Value = "2302"
$ws.Range("E18").Value = "2301"
$ws.Range("E19").Value = "2212"
$ws.Range("E20").Value = "2211"
$ws.Range("E21").Value = "2210"
$ws.Range("E22").Value = "2209"
$ws.Range("E23").Value = "2208"
$ws.Range("E24").Value = "2207"
$ws.Range("E25").Value = "2206"
